$d = $word.ActiveDocument

# The "Submitted On" date is typed into a floating text box (drawing
# canvas) that currently reads "/10/23" (the day was left blank).
# Locate that text box among the document's shapes and fill in the
# missing day ("26") at the very start of its text, giving "26/10/23".
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shape = $d.Shapes.Item($i)
    if ($shape.TextFrame.HasText) {
        $textRange = $shape.TextFrame.TextRange
        if ($textRange.Text -eq "/10/23") {
            $insertionPoint = $textRange.Duplicate
            $insertionPoint.Collapse(1)
            $insertionPoint.InsertBefore("26")
        }
    }
}
